# Remove the " Volume" run (accent3-colored) from the "Columns include: ..."
# text box and drop the now-trailing comma that used to separate
# "Adj Close" from "Volume", turning:
#   "...Adj Close, Volume, Name"
# into:
#   "...Adj Close, Name"
#
# The slide has two text boxes that happen to share this exact sentence;
# only the first one (shape Id=4, "文本框 3") is the intended target, so we
# stop at the first match encountered while walking the shapes collection.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }

    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    if ($full.IndexOf("Adj Close, Volume") -ge 0) {
        # 1) Delete the " Volume" substring (leading space + word), which is
        #    its own run with accent3 solid fill.
        $volPos0 = $full.IndexOf(" Volume")      # 0-based offset
        $volStart1 = $volPos0 + 1                # Characters() is 1-based
        $volLen = " Volume".Length

        $volRange = $tr.Characters($volStart1, $volLen)
        $volRange.Text = ""

        # 2) Delete the comma that used to precede " Volume", right after
        #    "Adj Close", so the run reads "Adj Close" with no trailing comma.
        $refreshed = $tr.Text
        $adjPos0 = $refreshed.IndexOf("Adj Close,")
        $commaStart1 = $adjPos0 + "Adj Close".Length + 1

        $commaRange = $tr.Characters($commaStart1, 1)
        $commaRange.Text = ""

        break
    }
}
